$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# Insert two new rows above the existing "StateBalancing_LA_GrossSalesTab" row
# (row 123) to hold the new LA R-1029-specific validation settings.
$ws.Rows("123:124").Insert()
$ws.Rows("123:124").RowHeight = 14.25

$ws.Range("A123").Value() = "StateBalancing_LA_GrossSalesTabR1029"
$ws.Range("A124").Value() = "StateBalancing_LA_GrossSalesClickOnTextR1029"

$ws.Range("C123").Value() = "Only for LA R-1029 and LA R-1029E"
$ws.Range("C124").Value() = "Only for LA R-1029 and LA R-1029E"

$ws.Range("B123").Value() = 1
$ws.Range("B124").Value() = "Gross sales of tangible###Gross Sales of Tangible###Gross sales of tangible###Gross sales of Tangible###GROSS SALES OF THE TANGIBLEGROSS SALES###Gross Sales###Gross sales###SERVICES AS REPORTED TO THE STATE###services as reported###RENTALS AND SERVICES AS REPORTED TO THE STATE OF LOUISIANA###and services as reported to the State of LA###SERVICES AS REPORTED TO THE STATE OF LOUISIANA###GROSS SALES OF TANGIBLE"

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 104
$ws.Range("B124").Select()
